$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("八位序列号收集收集结果yd5")

# --- Row 54: a "deleted" submission entry -------------------------------
# Reuse the strikethrough/gray formatting already used for the other
# deleted rows (e.g. row 13) instead of re-creating new style entries.
$ws.Range("A13:D13").Copy()
$ws.Range("A54:D54").PasteSpecial(-4122)
$ws.Cells.Item(54, 1).Value = "        "
$ws.Cells.Item(54, 2).Value = 45929.8563310185
$ws.Cells.Item(54, 3).Value = "已删除"
$ws.Cells.Item(54, 4).Value = "已删除"

# --- Row 55: a fresh submission ------------------------------------------
$ws.Cells.Item(55, 1).Value = "        "

# Reuse the date/time display format already used for the other rows.
$ws.Cells.Item(52, 2).Copy()
$ws.Cells.Item(55, 2).PasteSpecial(-4122)
$ws.Cells.Item(55, 2).Value = 45929.9812037037

$ws.Cells.Item(55, 3).Value = "eb0f7284"

# The QQ number looks numeric but must be stored as text (matches how the
# rest of the sheet keeps this column textual). Enter it as a formula that
# evaluates to a text string, then bake it down to a static value so it
# keeps its text type without picking up a new number-format style.
$ws.Cells.Item(55, 4).Formula = "=""781456741"""
$ws.Cells.Item(55, 4).Copy()
$ws.Cells.Item(55, 4).PasteSpecial(-4163)
